$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the data set. It lands at row 80 in the
# "Poroto granado" sheet, pushing the existing rows 80-95 down to 81-96.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new record's values.
$ws.Cells.Item(80, 1).Value = 7
$ws.Cells.Item(80, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(80, 3).Value = "Ñuble"
$ws.Cells.Item(80, 4).Value = 44641
$ws.Cells.Item(80, 5).Value = 16
$ws.Cells.Item(80, 6).Value = 100112030
$ws.Cells.Item(80, 7).Value = "Poroto granado"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 120
$ws.Cells.Item(80, 11).Value = 21000
$ws.Cells.Item(80, 12).Value = 22000
$ws.Cells.Item(80, 13).Value = 21500
$ws.Cells.Item(80, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(80, 16).Value = 860
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
